# Add "CUMPLIMIENTO MENSUAL" sheet (sheet3) with per-group budget-vs-sales
# compliance summary, as the last tab in the workbook.

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "CUMPLIMIENTO MENSUAL"

$ventasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# --- column widths (A:F) ---
$ws.Columns.Item(1).ColumnWidth = 27.17
$ws.Columns.Item(2).ColumnWidth = 21.17
$ws.Columns.Item(3).ColumnWidth = 21.17
$ws.Columns.Item(4).ColumnWidth = 11.17
$ws.Columns.Item(5).ColumnWidth = 21.17
$ws.Columns.Item(6).ColumnWidth = 27.17

# --- header row ---
$ws.Range("A1").Value = "ASESOR"
$ws.Range("B1").Value = "GRUPO"
$ws.Range("C1").Value = "PRESUPUESTO"
$ws.Range("D1").Value = "VENTA"
$ws.Range("E1").Value = "POR CUMPLIR"
$ws.Range("F1").Value = "CUMPLIMIENTO"

# match the header formatting already used on the other tabs
$ventasPorGrupo.Range("A1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# --- per-group rows: asesor / grupo / presupuesto / venta / por cumplir / cumplimiento ---
$asesor = "RIOS CARRION ANGEL BENIGNO"

$groups = @(
    @("240X120 PORCELANATO", 1041.16, 0, 1041.16, 0),
    @("240X80 PORCELANATO", 8668.91, 0, 8668.91, 0),
    @("FREGADEROS DE COCINA", 372.993863046034, 0, 372.993863046034, 0),
    @("GRANITO", 238.32, 0, 238.32, 0),
    @("GRIFERIAS", 106.82, 0, 106.82, 0),
    @("INODOROS", 800, 0, 800, 0),
    @("LAVABOS", 625, 0, 625, 0),
    @("LED", 300, 0, 300, 0),
    @("NO RESURTIBLES", 650.25, 0, 650.25, 0),
    @("OTROS", 0, 0, 0, 0),
    @("PANELES DECORATIVOS", 350, 0, 350, 0),
    @("PANELES PU", 230, 0, 230, 0),
    @("PANELES PVC", 483, 0, 483, 0),
    @("PIEDRA SINTERIZADA", 2501.01, 0, 2501.01, 0),
    @("PORCELANATO", 28209.84, -22.68, 28232.52, -0.0008039747832671153),
    @("PUERTAS DE SEGURIDAD", 342, 0, 342, 0),
    @("SAL SOLUBLE", 2300, 0, 2300, 0)
)

$row = 2
foreach ($g in $groups) {
    $ws.Cells.Item($row, 1).Value = $asesor
    $ws.Cells.Item($row, 2).Value = $g[0]
    $ws.Cells.Item($row, 3).Value = $g[1]
    $ws.Cells.Item($row, 4).Value = $g[2]
    $ws.Cells.Item($row, 5).Value = $g[3]
    $ws.Cells.Item($row, 6).Value = $g[4]
    $row = $row + 1
}

$lastDataRow = $row - 1
$totalRow = $lastDataRow + 1

# --- totals row ---
$ws.Cells.Item($totalRow, 2).Value = "TOTAL"
$ws.Cells.Item($totalRow, 2).HorizontalAlignment = -4152
$ws.Cells.Item($totalRow, 3).Value = 47219.30386304604
$ws.Cells.Item($totalRow, 4).Value = -22.68
$ws.Cells.Item($totalRow, 5).Value = 47241.98386304604
$ws.Cells.Item($totalRow, 6).Value = -0.0004803120364878872

# --- number formats: currency for PRESUPUESTO/VENTA/POR CUMPLIR, percent for CUMPLIMIENTO ---
$ventasPorGrupo.Range("C2").Copy()
$ws.Range("C2:E" + $totalRow).PasteSpecial(-4122)

$ws.Range("F2:F" + $totalRow).NumberFormat = "0.00%"

$ws.Range("A1").Select()
